# Fix Kd e Kc
# Updates the recomputed kc (C), ke (D), Eha (F) and kh (G) values (and
# a couple of floating-point-noise kd (B) values) to reflect the corrected
# kc factor, per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (h = 5)
$ws.Range("C2").Value = 0.6427876096865394
$ws.Range("D2").Value = 0.1398499385295467
$ws.Range("F2").Value = 3221.010519643452
$ws.Range("G2").Value = 0.02796998770590933

# Row 3 (h = 7)
$ws.Range("B3").Value = 0.1773629620793186
$ws.Range("C3").Value = 0.3420201433256688
$ws.Range("D3").Value = 0.06066170571103372
$ws.Range("F3").Value = 1397.154652259455
$ws.Range("G3").Value = 0.008665957958719103

# Row 4 (h = 11)
$ws.Range("C4").Value = 0.3420201433256687
$ws.Range("D4").Value = 0.06066170571103378
$ws.Range("F4").Value = 1397.154652259456
$ws.Range("G4").Value = 0.005514700519184889

# Row 5 (h = 13)
$ws.Range("C5").Value = 0.6427876096865394
$ws.Range("D5").Value = 0.1398499385295467
$ws.Range("F5").Value = 3221.010519643452
$ws.Range("G5").Value = 0.0107576875791959

# Row 6 (h = 17)
$ws.Range("B6").Value = 0.9597950805239367
$ws.Range("D6").Value = 0.9452136366029493
$ws.Range("F6").Value = 21770.07082606125
$ws.Range("G6").Value = 0.05560080215311466
